# "Generate Report for Handback"
#
# Refresh the handoff/handback timestamps recorded in the per-locale
# status report sheets for the "0f185259-..." source file now that a new
# handback round has been generated for it.
#
#   zh-cn: Correspond Handoff Datetime (H2)  -> 2016-08-16 10:47:56
#          Correspond Handback DateTime (K2) -> 2016-08-16 10:48:28
#   de-de: Correspond Handoff Datetime (H2)  -> 2016-08-16 10:48:04
#          Correspond Handback DateTime (K2) -> 2016-08-16 10:48:36

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 10:47:56"
$wsZhCn.Range("K2").Value = "2016-08-16 10:48:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 10:48:04"
$wsDeDe.Range("K2").Value = "2016-08-16 10:48:36"
